# Updates the cryptos price/volume table (Sheet1) with the latest
# scraped values, mirroring the GitHub Actions refresh commit.
#
# Column D holds price strings that often LOOK numeric (e.g. "83.50",
# "0.07740"); assigning such a string straight to Range.Value makes Excel
# coerce it to a real number and silently drop formatting (trailing
# zeros, leading pattern, etc). To keep the exact original text we mark
# the cell as Text ("@") before the assignment, then reset the style back
# to Normal afterwards so no stray number-format/style residue is left
# behind - only the cell value changes, just like the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.424.87'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.848.14'
$ws.Range("E3").Value = '  +0.34%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6260'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.64%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07672'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2916'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.57%  '

$ws.Range("E10").Value = '  +1.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07740'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = '1.842.53'
$ws.Range("E12").Value = '  -0.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.029'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("E14").Value = '  +3.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6800'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.83%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.171'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.87%  '

$ws.Range("D18").Value = '29.451.29'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.32%  '

$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("E21").Value = '  +0.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.400'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '

$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1373'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.395'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.62%  '

$ws.Range("E28").Value = '  +5.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.464'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05681'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.116'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.49%  '

$ws.Range("E32").Value = '  +0.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.840'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.73%  '

$ws.Range("E34").Value = '  +0.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7080'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.774'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("D38").Value = '1.228.04'
$ws.Range("E38").Value = '  -0.93%  '

$ws.Range("E39").Value = '  -0.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.540'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9098'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.17%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.90%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000120'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.14%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.138'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.63%  '

$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4015'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.980'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.03%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1146'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.49%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.669'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.29%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05715'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.12%  '
